$d = $word.ActiveDocument

# --- Change 1: Paragraph 1 ("Para el estado tensorial...") -------------
# Collapse the runs/proofErr-wrapped fragments (s_xx, s_yy, t_xy) into a
# single plain run by replacing the whole paragraph's text in one go.
# (Cleared first so the runtime doesn't treat an identical-text
# assignment as a no-op and skip re-merging the underlying runs.)
$p1 = $d.Paragraphs.Item(1).Range
[void]$p1.MoveEnd(1, -1)
$p1.Text = ""
$p1.Text = "Para el estado tensorial a 1cm obtener las tensiones que se están aplicando a partir de la matriz que forma los mapas, así podemos saber s_xx, s_yy y t_xy. Y después podemos hacer el círculo de Mohr."

# --- Change 2: insert new bullet right after "Igual para el valor máximo..." ---
$p2 = $d.Paragraphs.Item(2).Range
[void]$p2.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(3).Range
[void]$newPara.MoveEnd(1, -1)
$newPara.Text = "Convertir punto máximo de pixel a mm para saber su ubicación exacta"

# --- Change 3: "Buscar e implementar tensor de deformaciones..." -------
# Merge the trailing runs (" para hacer los mapas por " / "mi" / " mismo")
# into a single run by replacing the whole paragraph's text in one go.
$p6 = $d.Paragraphs.Item(6).Range
[void]$p6.MoveEnd(1, -1)
$p6.Text = ""
$p6.Text = "Buscar e implementar tensor de deformaciones para hacer los mapas por mi mismo"
